# Fruta / hortaliza, semanal
# The weekly refresh reshuffles the existing data rows (2-18) of the single
# worksheet: each row's "observation" fields (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Precio $/Kg,
# Kg/unidad) move to a different row, while the descriptive/id columns
# (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID, Producto,
# Categoria ID, Categoria, Variedad, Origen) stay put since they are the
# same for every row in this sheet. Row 15 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 18

# Snapshot current ("before") values of the columns that move between rows.
$colD = @{}
$colL = @{}
$colM = @{}
$colN = @{}
$colO = @{}
$colP = @{}
$colQ = @{}
$colS = @{}
$colT = @{}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $colD[$r] = $ws.Cells.Item($r, 4).Value2()
    $colL[$r] = $ws.Cells.Item($r, 12).Value()
    $colM[$r] = $ws.Cells.Item($r, 13).Value()
    $colN[$r] = $ws.Cells.Item($r, 14).Value()
    $colO[$r] = $ws.Cells.Item($r, 15).Value()
    $colP[$r] = $ws.Cells.Item($r, 16).Value()
    $colQ[$r] = $ws.Cells.Item($r, 17).Value()
    $colS[$r] = $ws.Cells.Item($r, 19).Value()
    $colT[$r] = $ws.Cells.Item($r, 20).Value()
}

# Target row -> source row (which row's data ends up at the target row).
$mapping = @{
    2  = 10
    3  = 11
    4  = 8
    5  = 16
    6  = 17
    7  = 18
    8  = 5
    9  = 2
    10 = 12
    11 = 13
    12 = 14
    13 = 7
    14 = 9
    15 = 15
    16 = 3
    17 = 4
    18 = 6
}

foreach ($target in ($mapping.Keys | Sort-Object)) {
    $source = $mapping[$target]

    $ws.Cells.Item($target, 4).Value = $colD[$source]
    $ws.Cells.Item($target, 12).Value = $colL[$source]
    $ws.Cells.Item($target, 13).Value = $colM[$source]
    $ws.Cells.Item($target, 14).Value = $colN[$source]
    $ws.Cells.Item($target, 15).Value = $colO[$source]
    $ws.Cells.Item($target, 16).Value = $colP[$source]
    $ws.Cells.Item($target, 17).Value = $colQ[$source]
    $ws.Cells.Item($target, 19).Value = $colS[$source]
    $ws.Cells.Item($target, 20).Value = $colT[$source]
}
